# Breadcrumbs.xlsx — "Chart" sheet GSC export refresh.
# The rolling 89-day window advances by two days: the two oldest dates
# (2025-11-08, 2025-11-09) drop off the front and two new dates
# (2026-02-05, 2026-02-06) are appended, so every row's Date (col A) and
# Items (col C) shift up by two and the last two rows become new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# New Date column values for rows 2..90 (col A), in order.
$newDates = @(
    "2025-11-10","2025-11-11","2025-11-12","2025-11-13","2025-11-14","2025-11-15","2025-11-16","2025-11-17","2025-11-18",
    "2025-11-19","2025-11-20","2025-11-21","2025-11-22","2025-11-23","2025-11-24","2025-11-25","2025-11-26","2025-11-27",
    "2025-11-28","2025-11-29","2025-11-30","2025-12-01","2025-12-02","2025-12-03","2025-12-04","2025-12-05","2025-12-06",
    "2025-12-07","2025-12-08","2025-12-09","2025-12-10","2025-12-11","2025-12-12","2025-12-13","2025-12-14","2025-12-15",
    "2025-12-16","2025-12-17","2025-12-18","2025-12-19","2025-12-20","2025-12-21","2025-12-22","2025-12-23","2025-12-24",
    "2025-12-25","2025-12-26","2025-12-27","2025-12-28","2025-12-29","2025-12-30","2025-12-31","2026-01-01","2026-01-02",
    "2026-01-03","2026-01-04","2026-01-05","2026-01-06","2026-01-07","2026-01-08","2026-01-09","2026-01-10","2026-01-11",
    "2026-01-12","2026-01-13","2026-01-14","2026-01-15","2026-01-16","2026-01-17","2026-01-18","2026-01-19","2026-01-20",
    "2026-01-21","2026-01-22","2026-01-23","2026-01-24","2026-01-25","2026-01-26","2026-01-27","2026-01-28","2026-01-29",
    "2026-01-30","2026-01-31","2026-02-01","2026-02-02","2026-02-03","2026-02-04","2026-02-05","2026-02-06"
)

# New Items column values for rows 2..90 (col C), in order.
$newItems = @(
    54.0,47.0,43.0,41.0,38.0,35.0,31.0,29.0,26.0,26.0,
    25.0,26.0,26.0,25.0,25.0,27.0,27.0,27.0,27.0,27.0,
    27.0,27.0,27.0,26.0,25.0,25.0,25.0,26.0,26.0,27.0,
    29.0,29.0,30.0,30.0,31.0,31.0,32.0,31.0,31.0,32.0,
    32.0,32.0,32.0,30.0,31.0,32.0,32.0,28.0,28.0,28.0,
    28.0,30.0,29.0,28.0,28.0,27.0,27.0,27.0,27.0,27.0,
    27.0,26.0,26.0,26.0,26.0,26.0,26.0,26.0,25.0,25.0,
    26.0,25.0,24.0,23.0,24.0,24.0,24.0,25.0,26.0,26.0,
    28.0,28.0,28.0,28.0,28.0,28.0,28.0,28.0,28.0
)

# Column A holds dates stored as literal text (e.g. "2025-11-10"), not real
# Excel date serials. Pre-formatting the range as Text keeps COM from
# auto-converting the assigned strings into date values.
$ws.Range("A2:A90").NumberFormat = "@"

for ($i = 0; $i -lt $newDates.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $newDates[$i]
    $ws.Cells.Item($row, 3).Value = $newItems[$i]
}
